# Apply the "release 1.6.3" update to CreateOrder-Event.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OrderCreated-Event")

# Update the MessageType value in row 3 (K3): JSONMessageType -> JSONType
$ws.Range("K3").Value = "JSONType"

# Reposition the view: scroll so column I row 3 is the top-left visible cell,
# and move the active selection to K4
$ws.Activate()
$ws.Range("K4").Select()
$excel.ActiveWindow.ScrollColumn = $ws.Range("I3").Column
$excel.ActiveWindow.ScrollRow = $ws.Range("I3").Row
